# Journal de travail - Luuk Mueller
# Adds two new weekly journal blocks (rover keyboard-control pages) to the
# "documentation/3_1_Journal-Luuk-Mueller.xlsx" journal sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Block 1 (rows 34-39): "Développement de la page pour conduire le rover
# au clavier" - single day entry on 08/01/2026.
# ---------------------------------------------------------------------
$ws.Range("A34").Value = 46030
$ws.Range("B34").Value = "Développement de la page pour conduire le rover au clavier"
$ws.Range("D34").Value = 3.5

# ---------------------------------------------------------------------
# Block 2 (rows 41-46): second week of the same feature - keyboard /
# controller choice pages, on 09/01/2026.
# ---------------------------------------------------------------------
$ws.Range("A41").Value = 46031
$ws.Range("B41").Value = "Dévélopper la page du clavier"
$ws.Range("D41").Value = 4

$ws.Range("B42").Value = "Fait la page pour choisir entre clavier et manette"
$ws.Range("D42").Value = 0.5

$ws.Range("B43").Value = "Merge entre dashboard et le reste"
$ws.Range("D43").Value = 0.5

$ws.Range("B44").Value = "Dévélopper la page de le manette"
$ws.Range("D44").Value = 2

# Reflection text for block 1 (row 40 header/merged cell).
$ws.Range("B40").Value = "Aujourd'hui, j'ai créé un site web permettant de contrôler le rover via le clavier. Cela m'a pris pas mal de temps, car j'ai dû attribuer un identifiant à chaque élément dans les fichiers HTML. Après quelques ajustements, tout fonctionne désormais correctement et le rover répond parfaitement aux commandes du clavier."
$ws.Rows.Item(40).RowHeight = 39

# Reflection text for block 2 (row 47 header/merged cell).
$ws.Range("B47").Value = "Aujourd'hui, j'ai bien avancé sur la création des pages de contrôle pour le rover. J'ai réussi à implémenter les deux options de contrôle, et l'intégration entre le tableau de bord et les autres pages s'est bien passée. Le plus gros défi a été de m'assurer que l'interface soit fluide et intuitive, surtout avec la gestion des choix entre clavier et manette. Tout fonctionne comme prévu, mais je dois encore tester certaines interactions pour être sûr de la stabilité."
$ws.Rows.Item(47).RowHeight = 50.25

# ---------------------------------------------------------------------
# Rows 35 and 36 of block 1 are left without any work-item text, so the
# previously merged/empty B:C cells are cleared out entirely (matching the
# rest of the un-used rows further down the sheet).
# ---------------------------------------------------------------------
$ws.Range("B35:C35").UnMerge()
$ws.Range("B35:C35").Clear()
$ws.Range("B36:C36").UnMerge()
$ws.Range("B36:C36").Clear()

# Row 44 re-uses the "first row of block" border style (same as row 41)
# rather than the plain repeated-row style, so copy formats across.
$ws.Range("B41:C41").Copy()
$ws.Range("B44:C44").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# Update the view so the window shows the newly added rows, mirroring
# where the author was working when they saved the file.
# ---------------------------------------------------------------------
$ws.Range("A48:A54").Select()
$excel.ActiveWindow.ScrollRow = 25
